$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.062.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "'1.862.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'311.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("D8").Value = "'0.3916"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "'0.09655"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +24.71%  "
$ws.Range("D10").Value = "'1.127"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("D11").Value = "'40.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "'6.445"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "'20.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").Value = "'1.867.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.94%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'7.359"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "'0.00001124"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.78%  "
$ws.Range("D18").Value = "'92.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "'0.06600"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "'17.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").Value = "'6.115"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").Value = "'28.144.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'11.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'2.285"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").Value = "'2.533"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.38%  "
$ws.Range("D27").Value = "'2.084.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("D28").Value = "'21.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.39%  "
$ws.Range("D29").Value = "'157.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").Value = "'127.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").Value = "'1.054"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").Value = "'5.598"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "'3.615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").Value = "'0.06734"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("D36").Value = "'9.439"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'0.02381"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").Value = "'0.2169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "'4.988"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").Value = "'0.6259"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").Value = "'1.003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "'13.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").Value = "'0.5978"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("D46").Value = "'3.662"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "'1.270"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").Value = "'124.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "'1.978"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("D50").Value = "'1.193"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "'0.06826"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.24%  "
